$wb = $excel.ActiveWorkbook

# --- Update "Metadata" sheet: refresh timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "29 Oct 2025, 08:39 AM"

# --- Update "1 Month Performance" sheet: refreshed % Change values (and reordered ties) ---
$perf = $wb.Worksheets.Item("1 Month Performance")

$perf.Range("C4").Value = 78.2013

$perf.Range("C6").Value = 66.7749

$perf.Range("C7").Value = 65.8165

$perf.Range("C10").Value = 52.479

$perf.Range("C12").Value = 45.9669

$perf.Range("C13").Value = 40.4222

$perf.Range("C14").Value = 40.3617

$perf.Range("C15").Value = 38.6669

$perf.Range("C17").Value = 38.0056

$perf.Range("C18").Value = 37.4816

$perf.Range("C19").Value = 36.8385

$perf.Range("B20").Value = "TVSELECT"
$perf.Range("C20").Value = 36.694

$perf.Range("B21").Value = "RAMAPHO"
$perf.Range("C21").Value = 36.6878

$perf.Range("C22").Value = 36.602

$perf.Range("C23").Value = 36.2847

$perf.Range("C24").Value = 36.2565

$perf.Range("C31").Value = 30.1764

$perf.Range("C34").Value = 27.8206

$perf.Range("C35").Value = 27.3801

$perf.Range("B36").Value = "ADANIPOWER"
$perf.Range("C36").Value = 27.2607

$perf.Range("B37").Value = "BHARATWIRE"
$perf.Range("C37").Value = 27.0597

$perf.Range("C38").Value = 26.8203

$perf.Range("B39").Value = "AVALON"
$perf.Range("C39").Value = 26.4978

$perf.Range("B40").Value = "HATSUN"
$perf.Range("C40").Value = 26.4313

$perf.Range("C41").Value = 25.7799

$perf.Range("C42").Value = 25.5353

$perf.Range("C43").Value = 25.2816

$perf.Range("C44").Value = 24.9687

$perf.Range("C45").Value = 24.8385

$perf.Range("C47").Value = 24.0369

$perf.Range("C49").Value = 23.4873

$perf.Range("C50").Value = 23.0862

$perf.Range("C51").Value = 23.0656

$perf.Range("B53").Value = "ETHOSLTD"
$perf.Range("C53").Value = 22.4308

$perf.Range("B54").Value = "INDIANB"
$perf.Range("C54").Value = 22.0659

$perf.Range("B55").Value = "GUJTHEM"
$perf.Range("C55").Value = 22.0441

$perf.Range("B56").Value = "ORBTEXP"
$perf.Range("C56").Value = 21.5577

$perf.Range("B57").Value = "PRIVISCL"
$perf.Range("C57").Value = 21.3912

$perf.Range("B58").Value = "TDPOWERSYS"
$perf.Range("C58").Value = 21.3285

$perf.Range("C60").Value = 20.4266

$perf.Range("C61").Value = 20.1474

$perf.Range("C63").Value = 19.7715

$perf.Range("B64").Value = "GRMOVER"
$perf.Range("C64").Value = 19.6558

$perf.Range("B65").Value = "ATL"
$perf.Range("C65").Value = 19.6524

$perf.Range("C66").Value = 19.5864

$perf.Range("B67").Value = "MANAKCOAT"
$perf.Range("C67").Value = 19.483

$perf.Range("B68").Value = "CEATLTD"
$perf.Range("C68").Value = 19.3871

$perf.Range("C70").Value = 19.1635

$perf.Range("B71").Value = "FEDERALBNK"
$perf.Range("C71").Value = 19.1522

$perf.Range("B75").Value = "TINNARUBR"
$perf.Range("C75").Value = 18.6773

$perf.Range("B76").Value = "M&MFIN"
$perf.Range("C76").Value = 18.6598
